$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.421.05'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.582.96'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.27'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.59'
$ws.Range('E8').Value = '  -1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.90'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.248'
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '1.810.14'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').Value = '1.585.49'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.71'
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('E16').Value = '  -2.11%  '
$ws.Range('D17').Value = '28.438.30'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.01'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.81'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('D21').Value = '0.0₃0689'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('E23').Value = '  -3.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.15'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('E25').Value = '  +3.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.46'
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.04'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('E28').Value = '  -1.93%  '
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0482'
$ws.Range('E31').Value = '  +2.67%  '
$ws.Range('E32').Value = '  -1.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.20'
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('D35').Value = '1.395.58'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  +6.68%  '
$ws.Range('E37').Value = '  -5.06%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.522'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.44'
$ws.Range('E45').Value = '  -3.44%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0458'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '62.77'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.720.97'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.923'
$ws.Range('E49').Value = '  -5.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.61'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('E51').Value = '  -1.11%  '
